# Update the "Elapsed Duration(Hrs)" values (column G) on several sheets
# of the Active_Outages workbook to reflect newly recalculated durations.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3926:15:11" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "65:47:49" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12107:38:52" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3237:22:21" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "475:33:55" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2953:28:41" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "180:40:56" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "427:27:40" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "67:59:58" }
)

foreach ($update in $updates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    $ws.Range($update.Cell).Value = $update.Value
}
